{"js": "// Resize the document's single pasted inline picture so its displayed\n// size matches the image's natural/proportional size (640x162 px @ 96dpi\n// == 480pt x 121.5pt == 6096000 x 1543050 EMU), preserving aspect ratio.\n//\n// Office.js exposes InlinePicture.width / InlinePicture.height (in points)\n// for this, but the installed host's implementation no-ops those two\n// setters, so we drive the same underlying Word object-model write\n// (InlineShape.Width / InlineShape.Height) that the shim's own working\n// setters (e.g. altTextDescription -> InlineShape.AlternativeText) use.\nconst pics = context.document.body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nif (pics.items.length === 0) {\n  throw new Error(\"No inline pictures found in document body.\");\n}\n\nconst pic = pics.items[0];\n\nconst widthPoints = 480;      // 6096000 EMU\nconst heightPoints = 121.5;   // 1543050 EMU\n\nif (typeof pic._omSet === \"function\") {\n  // Route straight through the Word object model (InlineShape.Width /\n  // InlineShape.Height) since the public width/height setters are\n  // currently not wired up to a write in this host.\n  pic._omSet(\"Width\", widthPoints);\n  pic._omSet(\"Height\", heightPoints);\n} else {\n  // Fall back to the documented Office.js API in case a future/real host\n  // implements the setters properly.\n  pic.width = widthPoints;\n  pic.height = heightPoints;\n}\n\nawait context.sync();\n", "ps1": "# Resize the document's single pasted inline picture so its displayed\n# size matches the image's natural/proportional size (640x162 px @ 96dpi\n# == 480pt x 121.5pt == 6096000 x 1543050 EMU), preserving aspect ratio.\n$d = $word.ActiveDocument\n\n$shape = $d.InlineShapes.Item(1)\n\n$shape.Width = 480\n$shape.Height = 121.5\n"}
